# Update cryptocurrency price/volume data per the source feed refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.878.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.80%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6881"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3047"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.081"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6813"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.441"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008291"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.892.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.076.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.11%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.462"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.58%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1475"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.792"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.545"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.208"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.149"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.181"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.838"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.138"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.217.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.699"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9379"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.740"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.66%  "

# Row 45
$ws.Range("E45").Value = "  -3.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5165"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "

# Row 47
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.976.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.19%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.496"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.96%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.746"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4185"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
